$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.173.99"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "2.480.62"
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.26%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "2.480.62"
$ws.Range("E9").Value = "  +0.57%  "

$ws.Range("E10").Value = "  +2.90%  "

$ws.Range("E11").Value = "  +1.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.94"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.333"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.42%  "

$ws.Range("D16").Value = "67.053.62"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000171"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").Value = "2.450.96"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "350.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.19%  "

$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.78%  "

$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.82"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.16%  "

$ws.Range("E27").Value = "  +2.49%  "

$ws.Range("D28").Value = "2.606.50"
$ws.Range("E28").Value = "  +0.81%  "

$ws.Range("E29").Value = "  +0.52%  "

$ws.Range("D30").Value = "0.0₃0910"
$ws.Range("E30").Value = "  +1.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "507.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.32%  "

$ws.Range("E33").Value = "  +0.98%  "

$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.29%  "

$ws.Range("E37").Value = "  -0.65%  "

$ws.Range("E38").Value = "  +0.65%  "

$ws.Range("E39").Value = "  -1.87%  "

$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("E43").Value = "  +1.39%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.89%  "

$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.20%  "

$ws.Range("D47").Value = "0.0₆0265"
$ws.Range("E47").Value = "  +4.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.34%  "
